$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-5
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 2
